# "Added support for metadata"
# Adds a new "process" column (F) with the full list of process types to the
# "Listen" sheet's dropdown-support table, and updates the view/selection
# state to match (Listen tab becomes the active tab).

$wb = $excel.ActiveWorkbook

$wsZuordnung = $wb.Worksheets.Item("Zuordnung")
$wsListen    = $wb.Worksheets.Item("Listen")

# --- Listen sheet: add the new "process" column (F) ------------------------

# Header (row 3, bold style already applied to the rest of the header row)
$wsListen.Cells.Item(3, 6).Value = "process"

# Values (rows 4-12)
$wsListen.Cells.Item(4, 6).Value  = "milling"
$wsListen.Cells.Item(5, 6).Value  = "drilling"
$wsListen.Cells.Item(6, 6).Value  = "grinding"
$wsListen.Cells.Item(7, 6).Value  = "turning"
$wsListen.Cells.Item(8, 6).Value  = "reaming"
$wsListen.Cells.Item(9, 6).Value  = "shaping"
$wsListen.Cells.Item(10, 6).Value = "thread_cutting"
$wsListen.Cells.Item(11, 6).Value = "thread_milling"
$wsListen.Cells.Item(12, 6).Value = "thread_forming"

# Column widths for the newly used columns E and F
$wsListen.Columns.Item(5).ColumnWidth = 13.2
$wsListen.Columns.Item(6).ColumnWidth = 23.9

# --- Selection / view state updates -----------------------------------------

# "Zuordnung": selection moves from P29 to G2:O2 (but it stays a non-active tab)
$wsZuordnung.Activate()
$wsZuordnung.Range("G2:O2").Select()

# "Listen": selection moves from D16 to F16, and becomes the active sheet/tab
# (this also leaves "Info" as no longer the active tab, matching the diff,
# since "Info"'s own selection (B7) is unchanged)
$wsListen.Activate()
$wsListen.Range("F16").Select()
